$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 203124
$ws.Cells.Item(2, 5).Value = 9979
$ws.Cells.Item(2, 6).Value = 9979
$ws.Cells.Item(2, 7).Value = -1374
$ws.Cells.Item(2, 8).Value = 332
$ws.Cells.Item(2, 9).Value = 653
$ws.Cells.Item(2, 10).Value = -321
$ws.Cells.Item(2, 11).Value = 313693
$ws.Cells.Item(2, 12).Value = 224670
$ws.Cells.Item(2, 13).Value = 89023
$ws.Cells.Item(2, 14).Value = 28894
$ws.Cells.Item(2, 15).Value = 60129
$ws.Cells.Item(2, 16).Value = 1348
$ws.Cells.Item(2, 17).Value = 6468
$ws.Cells.Item(2, 18).Value = -6110
$ws.Cells.Item(2, 19).Value = 4258
$ws.Cells.Item(2, 20).Value = 4533
$ws.Cells.Item(2, 21).Value = 1935
$ws.Cells.Item(2, 22).Value = 130405
$ws.Cells.Item(2, 23).Value = 4.91
$ws.Cells.Item(2, 24).Value = 0.16
$ws.Cells.Item(2, 25).Value = 2.21
$ws.Cells.Item(2, 26).Value = 0.11
$ws.Cells.Item(2, 27).Value = 252.37
$ws.Cells.Item(2, 28).Value = 2048.21
$ws.Cells.Item(2, 29).Value = 2478
$ws.Cells.Item(2, 30).Value = 41.78
$ws.Cells.Item(2, 31).Value = 140419
$ws.Cells.Item(2, 32).Value = 0.74
$ws.Cells.Item(2, 33).Value = 4000
$ws.Cells.Item(2, 34).Value = 3.86
$ws.Cells.Item(2, 35).Value = 126.59
$ws.Cells.Item(2, 36).Value = 21270888

# Row 3
$ws.Cells.Item(3, 4).Value = 169024
$ws.Cells.Item(3, 5).Value = 706
$ws.Cells.Item(3, 6).Value = 706
$ws.Cells.Item(3, 7).Value = -16278
$ws.Cells.Item(3, 8).Value = -17008
$ws.Cells.Item(3, 9).Value = -3912
$ws.Cells.Item(3, 10).Value = -13096
$ws.Cells.Item(3, 11).Value = 315563
$ws.Cells.Item(3, 12).Value = 231628
$ws.Cells.Item(3, 13).Value = 83935
$ws.Cells.Item(3, 14).Value = 25770
$ws.Cells.Item(3, 15).Value = 58165
$ws.Cells.Item(3, 16).Value = 1348
$ws.Cells.Item(3, 17).Value = -178
$ws.Cells.Item(3, 18).Value = -4458
$ws.Cells.Item(3, 19).Value = 11137
$ws.Cells.Item(3, 20).Value = 4075
$ws.Cells.Item(3, 21).Value = -4252
$ws.Cells.Item(3, 22).Value = 140561
$ws.Cells.Item(3, 23).Value = 0.42
$ws.Cells.Item(3, 24).Value = -10.06
$ws.Cells.Item(3, 25).Value = -14.31
$ws.Cells.Item(3, 26).Value = -5.41
$ws.Cells.Item(3, 27).Value = 275.96
$ws.Cells.Item(3, 28).Value = 1628.56
$ws.Cells.Item(3, 29).Value = -14671
$ws.Cells.Item(3, 30).Value = -6.03
$ws.Cells.Item(3, 31).Value = 128821
$ws.Cells.Item(3, 32).Value = 0.6899999999999999
$ws.Cells.Item(3, 33).Value = 4550
$ws.Cells.Item(3, 34).Value = 5.14
$ws.Cells.Item(3, 35).Value = -23.33
$ws.Cells.Item(3, 36).Value = 21270888

# Row 4
$ws.Cells.Item(4, 4).Value = 164703
$ws.Cells.Item(4, 5).Value = 9243
$ws.Cells.Item(4, 6).Value = 9243
$ws.Cells.Item(4, 7).Value = 17
$ws.Cells.Item(4, 8).Value = 504
$ws.Cells.Item(4, 9).Value = 1966
$ws.Cells.Item(4, 10).Value = -1462
$ws.Cells.Item(4, 11).Value = 286648
$ws.Cells.Item(4, 12).Value = 207639
$ws.Cells.Item(4, 13).Value = 79009
$ws.Cells.Item(4, 14).Value = 23915
$ws.Cells.Item(4, 15).Value = 55094
$ws.Cells.Item(4, 16).Value = 1348
$ws.Cells.Item(4, 17).Value = 9250
$ws.Cells.Item(4, 18).Value = 10612
$ws.Cells.Item(4, 19).Value = -22538
$ws.Cells.Item(4, 20).Value = 5043
$ws.Cells.Item(4, 21).Value = 4207
$ws.Cells.Item(4, 22).Value = 119843
$ws.Cells.Item(4, 23).Value = 5.61
$ws.Cells.Item(4, 24).Value = 0.31
$ws.Cells.Item(4, 25).Value = 7.91
$ws.Cells.Item(4, 26).Value = 0.17
$ws.Cells.Item(4, 27).Value = 262.81
$ws.Cells.Item(4, 28).Value = 1660.09
$ws.Cells.Item(4, 29).Value = 7384
$ws.Cells.Item(4, 30).Value = 14.22
$ws.Cells.Item(4, 31).Value = 121679
$ws.Cells.Item(4, 32).Value = 0.86
$ws.Cells.Item(4, 33).Value = 5100
$ws.Cells.Item(4, 34).Value = 4.86
$ws.Cells.Item(4, 35).Value = 51.09
$ws.Cells.Item(4, 36).Value = 20206888

# Row 5
$ws.Cells.Item(5, 4).Value = 169134
$ws.Cells.Item(5, 5).Value = 11676
$ws.Cells.Item(5, 6).Value = 11676
$ws.Cells.Item(5, 7).Value = 3559
$ws.Cells.Item(5, 8).Value = 459
$ws.Cells.Item(5, 9).Value = 406
$ws.Cells.Item(5, 10).Value = 54
$ws.Cells.Item(5, 11).Value = 287690
$ws.Cells.Item(5, 12).Value = 211656
$ws.Cells.Item(5, 13).Value = 76034
$ws.Cells.Item(5, 14).Value = 24741
$ws.Cells.Item(5, 15).Value = 51294
$ws.Cells.Item(5, 16).Value = 1348
$ws.Cells.Item(5, 17).Value = 6770
$ws.Cells.Item(5, 18).Value = -5560
$ws.Cells.Item(5, 19).Value = 2899
$ws.Cells.Item(5, 20).Value = 3948
$ws.Cells.Item(5, 21).Value = 2822
$ws.Cells.Item(5, 22).Value = 124275
$ws.Cells.Item(5, 23).Value = 6.9
$ws.Cells.Item(5, 24).Value = 0.27
$ws.Cells.Item(5, 25).Value = 1.67
$ws.Cells.Item(5, 26).Value = 0.16
$ws.Cells.Item(5, 27).Value = 278.37
$ws.Cells.Item(5, 28).Value = 1614.27
$ws.Cells.Item(5, 29).Value = 1594
$ws.Cells.Item(5, 30).Value = 70.89
$ws.Cells.Item(5, 31).Value = 125871
$ws.Cells.Item(5, 32).Value = 0.9
$ws.Cells.Item(5, 33).Value = 5100
$ws.Cells.Item(5, 34).Value = 4.51
$ws.Cells.Item(5, 35).Value = 259.3
$ws.Cells.Item(5, 36).Value = 19198003

# Row 6
$ws.Cells.Item(6, 4).Value = 181722
$ws.Cells.Item(6, 5).Value = 12159
$ws.Cells.Item(6, 6).Value = 12159
$ws.Cells.Item(6, 7).Value = 106
$ws.Cells.Item(6, 8).Value = -3405
$ws.Cells.Item(6, 9).Value = -1169
$ws.Cells.Item(6, 11).Value = 288803
$ws.Cells.Item(6, 12).Value = 217352
$ws.Cells.Item(6, 13).Value = 71450
$ws.Cells.Item(6, 14).Value = 20115
$ws.Cells.Item(6, 16).Value = 1348
$ws.Cells.Item(6, 17).Value = 10219
$ws.Cells.Item(6, 18).Value = -9431
$ws.Cells.Item(6, 19).Value = 108
$ws.Cells.Item(6, 20).Value = 3841
$ws.Cells.Item(6, 21).Value = 6378
$ws.Cells.Item(6, 22).Value = 122952
$ws.Cells.Item(6, 23).Value = 6.69
$ws.Cells.Item(6, 24).Value = -1.87
$ws.Cells.Item(6, 25).Value = -5.21
$ws.Cells.Item(6, 26).Value = -1.18
$ws.Cells.Item(6, 27).Value = 304.2
$ws.Cells.Item(6, 28).Value = 1276.69
$ws.Cells.Item(6, 29).Value = -4779
$ws.Cells.Item(6, 30).Value = -23.33
$ws.Cells.Item(6, 31).Value = 102338
$ws.Cells.Item(6, 33).Value = 5200
$ws.Cells.Item(6, 34).Value = 4.66
$ws.Cells.Item(6, 35).Value = -87.62
$ws.Cells.Item(6, 36).Value = 18238102

# Row 7
$ws.Cells.Item(7, 4).Value = 185682
$ws.Cells.Item(7, 5).Value = 12384
$ws.Cells.Item(7, 7).Value = 3432
$ws.Cells.Item(7, 8).Value = 1163
$ws.Cells.Item(7, 9).Value = 196
$ws.Cells.Item(7, 11).Value = 293206
$ws.Cells.Item(7, 12).Value = 222096
$ws.Cells.Item(7, 13).Value = 71110
$ws.Cells.Item(7, 14).Value = 18469
$ws.Cells.Item(7, 16).Value = 1287
$ws.Cells.Item(7, 17).Value = 1913
$ws.Cells.Item(7, 18).Value = -5977
$ws.Cells.Item(7, 19).Value = 4454
$ws.Cells.Item(7, 20).Value = 4202
$ws.Cells.Item(7, 21).Value = -1815
$ws.Cells.Item(7, 23).Value = 6.67
$ws.Cells.Item(7, 24).Value = 0.63
$ws.Cells.Item(7, 25).Value = 1.02
$ws.Cells.Item(7, 26).Value = 0.4
$ws.Cells.Item(7, 27).Value = 312.33
$ws.Cells.Item(7, 29).Value = 850
$ws.Cells.Item(7, 30).Value = 73.77
$ws.Cells.Item(7, 31).Value = 103709
$ws.Cells.Item(7, 32).Value = 0.6
$ws.Cells.Item(7, 33).Value = 5200
$ws.Cells.Item(7, 34).Value = 8.289999999999999
$ws.Cells.Item(7, 35).Value = 438.11

# Row 8
$ws.Cells.Item(8, 4).Value = 187062
$ws.Cells.Item(8, 5).Value = 13419
$ws.Cells.Item(8, 7).Value = 4426
$ws.Cells.Item(8, 8).Value = 2551
$ws.Cells.Item(8, 9).Value = 919
$ws.Cells.Item(8, 11).Value = 292857
$ws.Cells.Item(8, 12).Value = 220844
$ws.Cells.Item(8, 13).Value = 72014
$ws.Cells.Item(8, 14).Value = 18370
$ws.Cells.Item(8, 16).Value = 1287
$ws.Cells.Item(8, 17).Value = 8180
$ws.Cells.Item(8, 18).Value = -5116
$ws.Cells.Item(8, 19).Value = 806
$ws.Cells.Item(8, 20).Value = 4274
$ws.Cells.Item(8, 21).Value = 4782
$ws.Cells.Item(8, 23).Value = 7.17
$ws.Cells.Item(8, 24).Value = 1.36
$ws.Cells.Item(8, 25).Value = 4.99
$ws.Cells.Item(8, 26).Value = 0.87
$ws.Cells.Item(8, 27).Value = 306.67
$ws.Cells.Item(8, 29).Value = 4293
$ws.Cells.Item(8, 30).Value = 14.61
$ws.Cells.Item(8, 31).Value = 103155
$ws.Cells.Item(8, 32).Value = 0.61
$ws.Cells.Item(8, 33).Value = 5257
$ws.Cells.Item(8, 34).Value = 8.380000000000001
$ws.Cells.Item(8, 35).Value = 94.5

# Row 9
$ws.Cells.Item(9, 4).Value = 193944
$ws.Cells.Item(9, 5).Value = 14382
$ws.Cells.Item(9, 7).Value = 4985
$ws.Cells.Item(9, 8).Value = 2837
$ws.Cells.Item(9, 9).Value = 1181
$ws.Cells.Item(9, 11).Value = 295350
$ws.Cells.Item(9, 12).Value = 221853
$ws.Cells.Item(9, 13).Value = 73496
$ws.Cells.Item(9, 14).Value = 18591
$ws.Cells.Item(9, 16).Value = 1287
$ws.Cells.Item(9, 17).Value = 9364
$ws.Cells.Item(9, 18).Value = -6174
$ws.Cells.Item(9, 19).Value = 113
$ws.Cells.Item(9, 20).Value = 4344
$ws.Cells.Item(9, 21).Value = 4150
$ws.Cells.Item(9, 23).Value = 7.42
$ws.Cells.Item(9, 24).Value = 1.46
$ws.Cells.Item(9, 25).Value = 6.39
$ws.Cells.Item(9, 26).Value = 0.97
$ws.Cells.Item(9, 27).Value = 301.86
$ws.Cells.Item(9, 29).Value = 5515
$ws.Cells.Item(9, 30).Value = 11.37
$ws.Cells.Item(9, 31).Value = 104398
$ws.Cells.Item(9, 32).Value = 0.6
$ws.Cells.Item(9, 33).Value = 5314
$ws.Cells.Item(9, 34).Value = 8.48
$ws.Cells.Item(9, 35).Value = 74.34999999999999
